$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 2.027808622295978
$arr[0,1] = 9.960511441256548
$arr[0,2] = -7.953374226292953
$arr[0,3] = 0.2335353073468205
$arr[0,4] = 0.8732301461903136
$arr[0,5] = -0.9276358786007236
$arr[0,6] = -1.886570575107837
$arr[0,7] = 0.9471476497731522
$arr[0,8] = -0.7219570139330117
$arr[0,9] = 0.298215423810231
$ws.Range("B2:K2").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 7.710708601645155
$arr[0,1] = -10.20317706590435
$arr[0,2] = -2.016267532264572
$arr[0,3] = -1.376572693421079
$arr[0,4] = -3.177438718212116
$arr[0,5] = -4.136373414719229
$arr[0,6] = -1.30265518983824
$arr[0,7] = -2.971759853544404
$arr[0,8] = -1.951587415801161
$arr[0,9] = -3.665376872700473
$ws.Range("B3:K3").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -17.67232739676164
$arr[0,1] = -9.485417863121866
$arr[0,2] = -8.845723024278371
$arr[0,3] = -10.64658904906941
$arr[0,4] = -11.60552374557652
$arr[0,5] = -8.771805520695533
$arr[0,6] = -10.4409101844017
$arr[0,7] = -9.420737746658455
$arr[0,8] = -11.13452720355777
$arr[0,9] = -8.306171829958794
$ws.Range("B4:K4").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 8.685843937015603
$arr[0,1] = 9.325538775859098
$arr[0,2] = 7.52467275106806
$arr[0,3] = 6.565738054560947
$arr[0,4] = 9.399456279441935
$arr[0,5] = 7.730351615735771
$arr[0,6] = 8.750524053479014
$arr[0,7] = 7.036734596579703
$arr[0,8] = 9.865089970178675
$arr[0,9] = 7.629757370312555
$ws.Range("B5:K5").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1.185867228879948
$arr[0,1] = -0.6149987959110895
$arr[0,2] = -1.573933492418202
$arr[0,3] = 1.259784732462786
$arr[0,4] = -0.4093199312433776
$arr[0,5] = 0.6108525064998651
$arr[0,6] = -1.102936950399446
$arr[0,7] = 1.725418423199526
$arr[0,8] = -0.5099141766665937
$arr[0,9] = 0.8414208812510687
$ws.Range("B6:K6").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -2.762569267710482
$arr[0,1] = -3.721503964217595
$arr[0,2] = -0.8877857393366061
$arr[0,3] = -2.55689040304277
$arr[0,4] = -1.536717965299527
$arr[0,5] = -3.250507422198839
$arr[0,6] = -0.4221520485998669
$arr[0,7] = -2.657484648465986
$arr[0,8] = -1.306149590548324
$arr[0,9] = -1.927587520365226
$ws.Range("B7:K7").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.281473977694265
$arr[0,1] = 2.552244247186724
$arr[0,2] = 0.8831395834805599
$arr[0,3] = 1.903312021223803
$arr[0,4] = 0.1895225643244911
$arr[0,5] = 3.017877937923463
$arr[0,6] = 0.7825453380573438
$arr[0,7] = 2.133880395975006
$arr[0,8] = 1.512442466158104
$arr[0,9] = 2.013003163348936
$ws.Range("B8:K8").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 2.613418425600925
$arr[0,1] = 0.9443137618947609
$arr[0,2] = 1.964486199638004
$arr[0,3] = 0.2506967427386921
$arr[0,4] = 3.079052116337664
$arr[0,5] = 0.8437195164715449
$arr[0,6] = 2.195054574389207
$arr[0,7] = 1.573616644572305
$arr[0,8] = 2.074177341763137
$arr[0,9] = 1.406269269416995
$ws.Range("B9:K9").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -1.097142175261494
$arr[0,1] = -0.07696973751825081
$arr[0,2] = -1.790759194417562
$arr[0,3] = 1.03759617918141
$arr[0,4] = -1.19773642068471
$arr[0,5] = 0.1535986372329528
$arr[0,6] = -0.467839292583949
$arr[0,7] = 0.032721404606882
$arr[0,8] = -0.6351866677392595
$arr[0,9] = -0.2061856925012563
$ws.Range("B10:K10").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.5056247995153902
$arr[0,1] = -1.208164657383921
$arr[0,2] = 1.620190716215051
$arr[0,3] = -0.6151418836510686
$arr[0,4] = 0.7361931742665938
$arr[0,5] = 0.114755244449692
$arr[0,6] = 0.6153159416405229
$arr[0,7] = -0.05259213070561841
$arr[0,8] = 0.3764088445323847
$arr[0,9] = 0.4340803931105948
$ws.Range("B11:K11").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -1.303839698193279
$arr[0,1] = 1.524515675405693
$arr[0,2] = -0.7108169244604263
$arr[0,3] = 0.6405181334572361
$arr[0,4] = 0.01908020364033419
$arr[0,5] = 0.5196409008311652
$arr[0,6] = -0.1482671715149762
$arr[0,7] = 0.2807338037230269
$arr[0,8] = 0.338405352301237
$arr[0,9] = 0.4623717231395225
$ws.Range("B12:K12").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1.784934712404416
$arr[0,1] = -0.4503978874617036
$arr[0,2] = 0.9009371704559588
$arr[0,3] = 0.279499240639057
$arr[0,4] = 0.780059937829888
$arr[0,5] = 0.1121518654837466
$arr[0,6] = 0.5411528407217497
$arr[0,7] = 0.5988243892999598
$arr[0,8] = 0.7227907601382453
$arr[0,9] = -0.1706273630965465
$ws.Range("B13:K13").Value = $arr
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -1.285852883620195
$arr[0,1] = 0.06548217429746761
$arr[0,2] = -0.5559557555194342
$arr[0,3] = -0.0553950583286032
$arr[0,4] = -0.7233031306747446
$arr[0,5] = -0.2943021554367415
$arr[0,6] = -0.2366306068585314
$arr[0,7] = -0.1126642360202459
$arr[0,8] = -1.006082359255038
$arr[0,9] = -0.2359516323112753
$ws.Range("B14:K14").Value = $arr
$arr = New-Object 'object[,]' 1,9
$arr[0,0] = 0.1655615342000891
$arr[0,1] = -0.4558763956168127
$arr[0,2] = 0.04468430157401831
$arr[0,3] = -0.6232237707721231
$arr[0,4] = -0.19422279553412
$arr[0,5] = -0.1365512469559099
$arr[0,6] = -0.01258487611762438
$arr[0,7] = -0.9060029993524162
$arr[0,8] = -0.1358722724086538
$ws.Range("B15:J15").Value = $arr
$arr = New-Object 'object[,]' 1,8
$arr[0,0] = -0.1020898895371165
$arr[0,1] = 0.3984708076537146
$arr[0,2] = -0.2694372646924268
$arr[0,3] = 0.1595637105455762
$arr[0,4] = 0.2172352591237863
$arr[0,5] = 0.3412016299620719
$arr[0,6] = -0.55221649327272
$arr[0,7] = 0.2179142336710425
$ws.Range("B16:I16").Value = $arr
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.2502934172212692
$arr[0,1] = -0.4176146551248722
$arr[0,2] = 0.0113863201131309
$arr[0,3] = 0.06905786869134101
$arr[0,4] = 0.1930242395296265
$arr[0,5] = -0.7003938837051653
$arr[0,6] = 0.06973684323859711
$ws.Range("B17:H17").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = -0.5750606441290271
$arr[0,1] = -0.1460596688910241
$arr[0,2] = -0.08838812031281398
$arr[0,3] = 0.03557825052547153
$arr[0,4] = -0.8578398727093204
$arr[0,5] = -0.08770914576555788
$ws.Range("B18:G18").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0.3545997876350467
$arr[0,1] = 0.4122713362132568
$arr[0,2] = 0.5362377070515423
$arr[0,3] = -0.3571804161832495
$arr[0,4] = 0.4129503107605129
$ws.Range("B19:F19").Value = $arr
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 0.1319134556777877
$arr[0,1] = 0.2558798265160732
$arr[0,2] = -0.6375382967187186
$arr[0,3] = 0.1325924302250437
$ws.Range("B20:E20").Value = $arr
$arr = New-Object 'object[,]' 1,3
$arr[0,0] = 0.4278546843610848
$arr[0,1] = -0.465563438873707
$arr[0,2] = 0.3045672880700554
$ws.Range("B21:D21").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = -0.7714259786200386
$arr[0,1] = -0.001295251676276088
$ws.Range("B22:C22").Value = $arr
$ws.Range("B23").Value = 0.6110347010110101
$ws.Range("B24").Value = -0.343237405067616
